# Natmi following Dr Hou advice
# Update the sending/target cluster combinations for Ncam1-Ptprz1 so the
# sheet now reports all 9 combinations of the 3 clusters (ECs, FAPs, sCs)
# as both sender and target, instead of just 2 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value2 = "ECs"
$ws.Cells.Item(2,2).Value2 = "Ncam1"
$ws.Cells.Item(2,3).Value2 = "Ptprz1"
$ws.Cells.Item(2,4).Value2 = "ECs"
$ws.Cells.Item(2,5).Value2 = 2.0
$ws.Cells.Item(2,6).Value2 = 0.6666666666666666
$ws.Cells.Item(2,7).Value2 = 0.9404873333333333
$ws.Cells.Item(2,8).Value2 = 2.821462
$ws.Cells.Item(2,9).Value2 = 0.02000383747045655
$ws.Cells.Item(2,10).Value2 = 0.02000383747045654
$ws.Cells.Item(2,11).Value2 = 1.0
$ws.Cells.Item(2,12).Value2 = 0.3333333333333333
$ws.Cells.Item(2,13).Value2 = 0.03995766666666666
$ws.Cells.Item(2,14).Value2 = 0.119873
$ws.Cells.Item(2,15).Value2 = 0.005314930928687666
$ws.Cells.Item(2,16).Value2 = 0.005314930928687667
$ws.Cells.Item(2,17).Value2 = 0.03757967936955555
$ws.Cells.Item(2,18).Value2 = 0.338217114326
$ws.Cells.Item(2,19).Value2 = 0.0001063190144641707
$ws.Cells.Item(2,20).Value2 = 0.0001063190144641707

# Row 3
$ws.Cells.Item(3,1).Value2 = "ECs"
$ws.Cells.Item(3,2).Value2 = "Ncam1"
$ws.Cells.Item(3,3).Value2 = "Ptprz1"
$ws.Cells.Item(3,4).Value2 = "FAPs"
$ws.Cells.Item(3,5).Value2 = 2.0
$ws.Cells.Item(3,6).Value2 = 0.6666666666666666
$ws.Cells.Item(3,7).Value2 = 0.9404873333333333
$ws.Cells.Item(3,8).Value2 = 2.821462
$ws.Cells.Item(3,9).Value2 = 0.02000383747045655
$ws.Cells.Item(3,10).Value2 = 0.02000383747045654
$ws.Cells.Item(3,11).Value2 = 2.0
$ws.Cells.Item(3,12).Value2 = 0.6666666666666666
$ws.Cells.Item(3,13).Value2 = 0.03069133333333333
$ws.Cells.Item(3,14).Value2 = 0.092074
$ws.Cells.Item(3,15).Value2 = 0.004082378436578614
$ws.Cells.Item(3,16).Value2 = 0.004082378436578615
$ws.Cells.Item(3,17).Value2 = 0.02886481024311111
$ws.Cells.Item(3,18).Value2 = 0.259783292188
$ws.Cells.Item(3,19).Value2 = 0.00008166323473821508
$ws.Cells.Item(3,20).Value2 = 0.00008166323473821509

# Row 4
$ws.Cells.Item(4,1).Value2 = "ECs"
$ws.Cells.Item(4,2).Value2 = "Ncam1"
$ws.Cells.Item(4,3).Value2 = "Ptprz1"
$ws.Cells.Item(4,4).Value2 = "sCs"
$ws.Cells.Item(4,5).Value2 = 2.0
$ws.Cells.Item(4,6).Value2 = 0.6666666666666666
$ws.Cells.Item(4,7).Value2 = 0.9404873333333333
$ws.Cells.Item(4,8).Value2 = 2.821462
$ws.Cells.Item(4,9).Value2 = 0.02000383747045655
$ws.Cells.Item(4,10).Value2 = 0.02000383747045654
$ws.Cells.Item(4,11).Value2 = 3.0
$ws.Cells.Item(4,12).Value2 = 1.0
$ws.Cells.Item(4,13).Value2 = 7.447354000000001
$ws.Cells.Item(4,14).Value2 = 22.342062
$ws.Cells.Item(4,15).Value2 = 0.9906026906347337
$ws.Cells.Item(4,16).Value2 = 0.9906026906347338
$ws.Cells.Item(4,17).Value2 = 7.004142103849334
$ws.Cells.Item(4,18).Value2 = 63.03727893464401
$ws.Cells.Item(4,19).Value2 = 0.01981585522125416
$ws.Cells.Item(4,20).Value2 = 0.01981585522125416

# Row 5
$ws.Cells.Item(5,1).Value2 = "FAPs"
$ws.Cells.Item(5,2).Value2 = "Ncam1"
$ws.Cells.Item(5,3).Value2 = "Ptprz1"
$ws.Cells.Item(5,4).Value2 = "ECs"
$ws.Cells.Item(5,5).Value2 = 3.0
$ws.Cells.Item(5,6).Value2 = 1.0
$ws.Cells.Item(5,7).Value2 = 1.392600333333333
$ws.Cells.Item(5,8).Value2 = 4.177801000000001
$ws.Cells.Item(5,9).Value2 = 0.0296201232509638
$ws.Cells.Item(5,10).Value2 = 0.0296201232509638
$ws.Cells.Item(5,11).Value2 = 1.0
$ws.Cells.Item(5,12).Value2 = 0.3333333333333333
$ws.Cells.Item(5,13).Value2 = 0.03995766666666666
$ws.Cells.Item(5,14).Value2 = 0.119873
$ws.Cells.Item(5,15).Value2 = 0.005314930928687666
$ws.Cells.Item(5,16).Value2 = 0.005314930928687667
$ws.Cells.Item(5,17).Value2 = 0.05564505991922222
$ws.Cells.Item(5,18).Value2 = 0.5008055392730001
$ws.Cells.Item(5,19).Value2 = 0.0001574289091780882
$ws.Cells.Item(5,20).Value2 = 0.0001574289091780882

# Row 6
$ws.Cells.Item(6,1).Value2 = "FAPs"
$ws.Cells.Item(6,2).Value2 = "Ncam1"
$ws.Cells.Item(6,3).Value2 = "Ptprz1"
$ws.Cells.Item(6,4).Value2 = "FAPs"
$ws.Cells.Item(6,5).Value2 = 3.0
$ws.Cells.Item(6,6).Value2 = 1.0
$ws.Cells.Item(6,7).Value2 = 1.392600333333333
$ws.Cells.Item(6,8).Value2 = 4.177801000000001
$ws.Cells.Item(6,9).Value2 = 0.0296201232509638
$ws.Cells.Item(6,10).Value2 = 0.0296201232509638
$ws.Cells.Item(6,11).Value2 = 2.0
$ws.Cells.Item(6,12).Value2 = 0.6666666666666666
$ws.Cells.Item(6,13).Value2 = 0.03069133333333333
$ws.Cells.Item(6,14).Value2 = 0.092074
$ws.Cells.Item(6,15).Value2 = 0.004082378436578614
$ws.Cells.Item(6,16).Value2 = 0.004082378436578615
$ws.Cells.Item(6,17).Value2 = 0.04274076103044445
$ws.Cells.Item(6,18).Value2 = 0.3846668492740001
$ws.Cells.Item(6,19).Value2 = 0.0001209205524485355
$ws.Cells.Item(6,20).Value2 = 0.0001209205524485355

# Row 7
$ws.Cells.Item(7,1).Value2 = "FAPs"
$ws.Cells.Item(7,2).Value2 = "Ncam1"
$ws.Cells.Item(7,3).Value2 = "Ptprz1"
$ws.Cells.Item(7,4).Value2 = "sCs"
$ws.Cells.Item(7,5).Value2 = 3.0
$ws.Cells.Item(7,6).Value2 = 1.0
$ws.Cells.Item(7,7).Value2 = 1.392600333333333
$ws.Cells.Item(7,8).Value2 = 4.177801000000001
$ws.Cells.Item(7,9).Value2 = 0.0296201232509638
$ws.Cells.Item(7,10).Value2 = 0.0296201232509638
$ws.Cells.Item(7,11).Value2 = 3.0
$ws.Cells.Item(7,12).Value2 = 1.0
$ws.Cells.Item(7,13).Value2 = 7.447354000000001
$ws.Cells.Item(7,14).Value2 = 22.342062
$ws.Cells.Item(7,15).Value2 = 0.9906026906347337
$ws.Cells.Item(7,16).Value2 = 0.9906026906347338
$ws.Cells.Item(7,17).Value2 = 10.37118766285134
$ws.Cells.Item(7,18).Value2 = 93.34068896566203
$ws.Cells.Item(7,19).Value2 = 0.02934177378933718
$ws.Cells.Item(7,20).Value2 = 0.02934177378933718

# Row 8
$ws.Cells.Item(8,1).Value2 = "sCs"
$ws.Cells.Item(8,2).Value2 = "Ncam1"
$ws.Cells.Item(8,3).Value2 = "Ptprz1"
$ws.Cells.Item(8,4).Value2 = "ECs"
$ws.Cells.Item(8,5).Value2 = 3.0
$ws.Cells.Item(8,6).Value2 = 1.0
$ws.Cells.Item(8,7).Value2 = 44.682258
$ws.Cells.Item(8,8).Value2 = 134.046774
$ws.Cells.Item(8,9).Value2 = 0.9503760392785797
$ws.Cells.Item(8,10).Value2 = 0.9503760392785796
$ws.Cells.Item(8,11).Value2 = 1.0
$ws.Cells.Item(8,12).Value2 = 0.3333333333333333
$ws.Cells.Item(8,13).Value2 = 0.03995766666666666
$ws.Cells.Item(8,14).Value2 = 0.119873
$ws.Cells.Item(8,15).Value2 = 0.005314930928687666
$ws.Cells.Item(8,16).Value2 = 0.005314930928687667
$ws.Cells.Item(8,17).Value2 = 1.785398771078
$ws.Cells.Item(8,18).Value2 = 16.068588939702
$ws.Cells.Item(8,19).Value2 = 0.005051183005045408
$ws.Cells.Item(8,20).Value2 = 0.005051183005045408

# Row 9
$ws.Cells.Item(9,1).Value2 = "sCs"
$ws.Cells.Item(9,2).Value2 = "Ncam1"
$ws.Cells.Item(9,3).Value2 = "Ptprz1"
$ws.Cells.Item(9,4).Value2 = "FAPs"
$ws.Cells.Item(9,5).Value2 = 3.0
$ws.Cells.Item(9,6).Value2 = 1.0
$ws.Cells.Item(9,7).Value2 = 44.682258
$ws.Cells.Item(9,8).Value2 = 134.046774
$ws.Cells.Item(9,9).Value2 = 0.9503760392785797
$ws.Cells.Item(9,10).Value2 = 0.9503760392785796
$ws.Cells.Item(9,11).Value2 = 2.0
$ws.Cells.Item(9,12).Value2 = 0.6666666666666666
$ws.Cells.Item(9,13).Value2 = 0.03069133333333333
$ws.Cells.Item(9,14).Value2 = 0.092074
$ws.Cells.Item(9,15).Value2 = 0.004082378436578614
$ws.Cells.Item(9,16).Value2 = 0.004082378436578615
$ws.Cells.Item(9,17).Value2 = 1.371358074364
$ws.Cells.Item(9,18).Value2 = 12.342222669276
$ws.Cells.Item(9,19).Value2 = 0.003879794649391863
$ws.Cells.Item(9,20).Value2 = 0.003879794649391864

# Row 10
$ws.Cells.Item(10,1).Value2 = "sCs"
$ws.Cells.Item(10,2).Value2 = "Ncam1"
$ws.Cells.Item(10,3).Value2 = "Ptprz1"
$ws.Cells.Item(10,4).Value2 = "sCs"
$ws.Cells.Item(10,5).Value2 = 3.0
$ws.Cells.Item(10,6).Value2 = 1.0
$ws.Cells.Item(10,7).Value2 = 44.682258
$ws.Cells.Item(10,8).Value2 = 134.046774
$ws.Cells.Item(10,9).Value2 = 0.9503760392785797
$ws.Cells.Item(10,10).Value2 = 0.9503760392785796
$ws.Cells.Item(10,11).Value2 = 3.0
$ws.Cells.Item(10,12).Value2 = 1.0
$ws.Cells.Item(10,13).Value2 = 7.447354000000001
$ws.Cells.Item(10,14).Value2 = 22.342062
$ws.Cells.Item(10,15).Value2 = 0.9906026906347337
$ws.Cells.Item(10,16).Value2 = 0.9906026906347338
$ws.Cells.Item(10,17).Value2 = 332.764592845332
$ws.Cells.Item(10,18).Value2 = 2994.881335607988
$ws.Cells.Item(10,19).Value2 = 0.9414450616241424
$ws.Cells.Item(10,20).Value2 = 0.9414450616241424
